$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Word"
$ws.Range("B1").Value = "Meaning"

$ws.Range("A2").Value = "('Word', ""('Word', 'Hello')"")"
$ws.Range("B2").Value = "('Meaning', ""('Meaning', 'A greeting')"")"

$ws.Range("A3").Value = "('Word', ""('Word', 'Goodbye')"")"
$ws.Range("B3").Value = "('Meaning', ""('Meaning', 'A farewell')"")"

$ws.Range("A4").Value = "('Word', ""('Word', 'Thank you')"")"
$ws.Range("B4").Value = "('Meaning', ""('Meaning', 'Gratitude')"")"

$ws.Range("A5").Value = "('Word', ""('Word', 'Dancer')"")"
$ws.Range("B5").Value = "('Meaning', ""('Meaning', 'dsf')"")"

$ws.Range("A6").Value = "('Word', 'dsfa')"
$ws.Range("B6").Value = "('Meaning', 'adsfasf')"
